# Generate Report for Handback
#
# The hand-off for 49906aa1-7bfa-4bd9-ba7c-060d782057c6.md failed its
# handback transform because the handback file name didn't match the
# handoff file name. Update the status for that row (Overview + the two
# per-locale sheets) and record the error detail on each locale sheet.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"
$zhCnError  = "Handback file name: d2dxaph2.0sx is different with handoff file name: 49906aa1-7bfa-4bd9-ba7c-060d782057c6.79a23d1e2efe6af140fe58c6d3a1a9f105ff0667.zh-cn."
$deDeError  = "Handback file name: d2dxaph2.0sx is different with handoff file name: 49906aa1-7bfa-4bd9-ba7c-060d782057c6.79a23d1e2efe6af140fe58c6d3a1a9f105ff0667.de-de."

# Overview sheet: row for 49906aa1-7bfa-4bd9-ba7c-060d782057c6.md is row 8.
# Column E = zh-cn status, column F = de-de status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E8").Value = $statusText
$wsOverview.Range("F8").Value = $statusText

# zh-cn sheet: row 8 is the same file. Column C = Status, column R = Error Detail.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C8").Value = $statusText
$wsZhCn.Range("R8").Value = $zhCnError

# de-de sheet: row 8 is the same file. Column C = Status, column R = Error Detail.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C8").Value = $statusText
$wsDeDe.Range("R8").Value = $deDeError
